# The file "b814fa4b-2a4a-4bb4-a9c7-189647d9a0e4" has finished translation and
# is now ready for handoff: update its Status cell on each sheet and refresh
# the "Latest Handoff Datetime" stamps on the per-locale sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B6").Value = "Ready for handoff"
$overview.Range("C6").Value = "Ready for handoff"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B6").Value = "Ready for handoff"
$zhcn.Range("D6").Value = "2016-03-09 09:28:16"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B6").Value = "Ready for handoff"
$dede.Range("D6").Value = "2016-03-09 09:28:24"
